$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K header, copying format from J1
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("K1").Value = "Buckle ratio"

# Update changed cells in columns C, F, G (rows 2-21)
$ws.Range("C2").Value = 337.2695627276308
$ws.Range("F2").Value = 175.0122316608695
$ws.Range("G2").Value = [double]"-4.293176737148315e-11"
$ws.Range("C3").Value = 625.8612525617799
$ws.Range("F3").Value = 397.1059129020874
$ws.Range("G3").Value = 245.6123408854625
$ws.Range("C4").Value = 822.9134908162224
$ws.Range("F4").Value = 648.3402844041186
$ws.Range("G4").Value = 802.5533206764782
$ws.Range("C5").Value = 897.9960599335344
$ws.Range("F5").Value = 897.7637848035575
$ws.Range("G5").Value = 1714.34337620186
$ws.Range("C6").Value = 926.9347789619342
$ws.Range("F6").Value = 1156.953952186281
$ws.Range("G6").Value = 2978.324550276093
$ws.Range("C7").Value = 935.5258614047922
$ws.Range("F7").Value = 1425.831136760593
$ws.Range("G7").Value = 4607.22523950447
$ws.Range("C8").Value = 930.9463650052438
$ws.Range("F8").Value = 1701.377148470604
$ws.Range("G8").Value = 6614.683984933038
$ws.Range("C9").Value = 872.831072611623
$ws.Range("F9").Value = 1987.147194106904
$ws.Range("G9").Value = 9007.839350409757
$ws.Range("C10").Value = 903.6782435342604
$ws.Range("F10").Value = 2297.737444150311
$ws.Range("G10").Value = 11649.68579295559
$ws.Range("C11").Value = 853.1687190743983
$ws.Range("F11").Value = 2614.389296532579
$ws.Range("G11").Value = 14983.81667906232
$ws.Range("C12").Value = 853.417191681895
$ws.Range("F12").Value = -2614.439928077591
$ws.Range("G12").Value = 18646.69199201048
$ws.Range("C13").Value = 904.2670374911922
$ws.Range("F13").Value = -2297.963454271599
$ws.Range("G13").Value = 15077.63455762345
$ws.Range("C14").Value = 873.5564153245804
$ws.Range("F14").Value = -1987.510904179655
$ws.Range("G14").Value = 11678.79904100214
$ws.Range("C15").Value = 931.6159064490627
$ws.Range("F15").Value = -1701.676930506061
$ws.Range("G15").Value = 9012.017504124104
$ws.Range("C16").Value = 936.2359146908516
$ws.Range("F16").Value = -1426.13868462064
$ws.Range("G16").Value = 6616.189777222519
$ws.Range("C17").Value = 927.6396865498834
$ws.Range("F17").Value = -1157.204843849535
$ws.Range("G17").Value = 4608.298301605118
$ws.Range("C18").Value = 898.692459366369
$ws.Range("F18").Value = -897.9716041541654
$ws.Range("G18").Value = 2979.044593352113
$ws.Range("C19").Value = 823.5477617982294
$ws.Range("F19").Value = -648.4794598773691
$ws.Range("G19").Value = 1715.027634685571
$ws.Range("C20").Value = 626.3769706185507
$ws.Range("F20").Value = -397.215658555278
$ws.Range("G20").Value = 803.1070853604114
$ws.Range("C21").Value = 337.5478961335763
$ws.Range("F21").Value = -175.0694776767505
$ws.Range("G21").Value = 245.6809126371058

# New column K values (Buckle ratio)
$ws.Range("K2").Value = 0.9503561321756911
$ws.Range("K3").Value = 4.047704900001258
$ws.Range("K4").Value = 9.704344860264268
$ws.Range("K5").Value = 18.04893093899456
$ws.Range("K6").Value = 29.11339035925919
$ws.Range("K7").Value = 43.00304931844122
$ws.Range("K8").Value = 59.80142576705639
$ws.Range("K9").Value = 79.43819985057675
$ws.Range("K10").Value = 103.2178955249094
$ws.Range("K11").Value = 131.5834403192706
$ws.Range("K12").Value = 131.6411232115622
$ws.Range("K13").Value = 103.3026518664712
$ws.Range("K14").Value = 79.50517434511474
$ws.Range("K15").Value = 59.84782283564201
$ws.Range("K16").Value = 43.0368592634219
$ws.Range("K17").Value = 29.13648317783106
$ws.Range("K18").Value = 18.06340638119294
$ws.Range("K19").Value = 9.712277375218621
$ws.Range("K20").Value = 4.051097989488993
$ws.Range("K21").Value = 0.9511505493973414
